$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.101
$ws.Range("E2").Value = -0.07780000000000001
$ws.Range("G2").Value = 0.252101731968866
$ws.Range("H2").Value = 0.2521013550441944
$ws.Range("I2").Value = 0.1972852000527694
$ws.Range("J2").Value = 0.1405497290206908
$ws.Range("K2").Value = 555.03
$ws.Range("L2").Value = 0.05230112512014473
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 2260.35
$ws.Range("V2").Value = 0.3004302398683628
$ws.Range("W2").Value = 0.1437819776547979
$ws.Range("X2").Value = 0.07635137176111639
$ws.Range("Y2").Value = 0.06743060589368155
$ws.Range("Z2").Value = 1.435978901845954
$ws.Range("AA2").Value = 0.1832507642036471
$ws.Range("AB2").Value = 0.05004852750165143
$ws.Range("AC2").Value = 0.1352663844618457
$ws.Range("AD2").Value = 7691.29
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 7691.29
$ws.Range("AG2").Value = 5430.940000000001
$ws.Range("AH2").Value = 0.5055070653959908
$ws.Range("AI2").Value = 0.4125544638791102
$ws.Range("AJ2").Value = 0.4192270729043239
$ws.Range("AK2").Value = 0.331503765548376
$ws.Range("AL2").Value = 402.127
$ws.Range("AM2").Value = 401.902
$ws.Range("AN2").Value = 3.397107863272779
$ws.Range("AO2").Value = 5.206390021062997
$ws.Range("AP2").Value = 2.398750922012129
$ws.Range("AQ2").Value = 5.209304755885763
$ws.Range("D3").Value = 0.101
$ws.Range("E3").Value = -0.07780000000000001
$ws.Range("G3").Value = 0.16078125
$ws.Range("H3").Value = 0.1606770833333333
$ws.Range("I3").Value = 0.09531250000000001
$ws.Range("J3").Value = 0.06833726415094341
$ws.Range("K3").Value = 2.73
$ws.Range("L3").Value = 0.07109375
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 1.68
$ws.Range("V3").Value = 0.04375
$ws.Range("W3").Value = 0.21328125
$ws.Range("X3").Value = 0.05101576118657218
$ws.Range("Y3").Value = 0.1622654888134278
$ws.Range("Z3").Value = 2.681564245810056
$ws.Range("AA3").Value = 0.1832507642036471
$ws.Range("AB3").Value = 0.04749335870246536
$ws.Range("AC3").Value = 0.1357574055011818
$ws.Range("AD3").Value = 6.59
$ws.Range("AF3").Value = 6.59
$ws.Range("AG3").Value = 4.91
$ws.Range("AH3").Value = 0.1464769948877528
$ws.Range("AI3").Value = 0.2917220008853475
$ws.Range("AJ3").Value = 0.1133687370122374
$ws.Range("AK3").Value = 0.2348158775705404
$ws.Range("AL3").Value = 0.201
$ws.Range("AM3").Value = 0.177
$ws.Range("AN3").Value = 1.651629072681704
$ws.Range("AO3").Value = 18.2089552238806
$ws.Range("AP3").Value = 1.230576441102757
$ws.Range("AQ3").Value = 20.67796610169491
$ws.Range("D4").Value = 0.154
$ws.Range("E4").Value = -0.0101
$ws.Range("G4").Value = 0.2524387702915159
$ws.Range("H4").Value = 0.2524387702915159
$ws.Range("I4").Value = 0.1976966643304984
$ws.Range("J4").Value = 0.1292163711428461
$ws.Range("K4").Value = 550.8
$ws.Range("L4").Value = 0.05216600685696967
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 2256.1
$ws.Range("V4").Value = 0.3016458759509579
$ws.Range("W4").Value = 0.1437819776547979
$ws.Range("X4").Value = 0.07635137176111639
$ws.Range("Y4").Value = 0.06743060589368155
$ws.Range("Z4").Value = 1.434144221235212
$ws.Range("AA4").Value = 0.1853149119634971
$ws.Range("AB4").Value = 0.05004852750165143
$ws.Range("AC4").Value = 0.1352663844618457
$ws.Range("AD4").Value = 7671.1
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 7671.1
$ws.Range("AG4").Value = 5415
$ws.Range("AH4").Value = 0.5063298658781286
$ws.Range("AI4").Value = 0.4125578143487146
$ws.Range("AJ4").Value = 0.4199530024894721
$ws.Range("AK4").Value = 0.3314379449011195
$ws.Range("AL4").Value = 401.1
$ws.Range("AM4").Value = 401.1
$ws.Range("AN4").Value = 3.398653138983652
$ws.Range("AO4").Value = 5.204188481675392
$ws.Range("AP4").Value = 2.399096185370608
$ws.Range("AQ4").Value = 5.204188481675392
$ws.Range("D5").Value = -0.0767
$ws.Range("E5").Value = -0.245
$ws.Range("G5").Value = 0.2486842105263158
$ws.Range("H5").Value = 0.2486842105263158
$ws.Range("I5").Value = 0.1690789473684211
$ws.Range("J5").Value = 0.1296271929824561
$ws.Range("K5").Value = 1.5
$ws.Range("L5").Value = 0.09868421052631579
$ws.Range("U5").Value = 2.57
$ws.Range("V5").Value = 0.4276206322795341
$ws.Range("W5").Value = 0.1271186440677966
$ws.Range("X5").Value = 0.11305561922871
$ws.Range("Y5").Value = 0.01406302483908664
$ws.Range("Z5").Value = 1.117647058823529
$ws.Range("AA5").Value = 0.1448774509803921
$ws.Range("AB5").Value = 0.05570014211773108
$ws.Range("AC5").Value = 0.08917730886266106
$ws.Range("AD5").Value = 13.6
$ws.Range("AF5").Value = 13.6
$ws.Range("AG5").Value = 11.03
$ws.Range("AH5").Value = 0.6935237123916369
$ws.Range("AI5").Value = 0.5132075471698113
$ws.Range("AJ5").Value = 0.6473004694835681
$ws.Range("AK5").Value = 0.4609277058086084
$ws.Range("AL5").Value = 0.826
$ws.Range("AM5").Value = 0.625
$ws.Range("AN5").Value = 4.563758389261745
$ws.Range("AO5").Value = 3.11138014527845
$ws.Range("AP5").Value = 3.701342281879195
$ws.Range("AQ5").Value = 4.112

# Remove the buybacks_cash_returned (T) values for rows 2-4 entirely (cell removed from row)
$ws.Range("T2").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("T4").ClearContents()
